# Insert a new weekly data row at row 10 (pushes existing rows 10-121 down to 11-122)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the new week's record.
# Columns that repeat the same market/category/unit/origin info are copied
# from the neighbouring rows; D/J/K/L/M/P carry the new data point.
$ws.Range("A10").Value = 6
$ws.Range("B10").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C10").Value = "Metropolitana"
$ws.Range("D10").Value = 45022
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = 100114007
$ws.Range("G10").Value = "Jengibre"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 180
$ws.Range("K10").Value = 15000
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 15444
$ws.Range("N10").Value = "`$/caja 13 kilos"
$ws.Range("O10").Value = "Perú"
$ws.Range("P10").Value = 1188
$ws.Range("Q10").Value = 13
$ws.Range("R10").Value = "Hortaliza"
